$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by one day
# (45180 -> 45181) for every data row (rows 2 through 221).
$ws.Range("C2:C221").Value = 45181
